# Update "Project Review Spreadsheet" per latest team review edits.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Documentation section: expand project-log note to mention the Gantt chart.
$ws.Range("D31").Value = "Project log is updated every meeting - changes reflected in SPMP and Gantt chart."

# Security section updates.
$ws.Range("D35").Value = "Validate security and encryption by testing login / log off  feature."
$ws.Range("D36").Value = "Periodically update packages"
$ws.Range("D37").Value = "Preventing user from leaving website to ensure smooth client experience."

# Select the last-edited cell to mirror the saved view state.
$ws.Range("D36").Select()
